# ---------------------------------------------------------------------------
# feat: add 2022-Q4 data
#
# The workbook tracks a fund's quarterly shareholder snapshots. A new
# quarter ("2022-Q4") is inserted as a sheet right after the "总计"
# (summary) sheet and before the existing "2022-Q3" sheet, and the
# "总计" summary sheet gets a new leading row for the 2022-Q4 totals.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q4" sheet immediately before "2022-Q3", so tab
#    order becomes: 总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$ws = $wb.Worksheets.Add($q3)
$ws.Name = "2022-Q4"

# Match the page-margin convention used by the other quarter sheets
# (0.75in/0.75in/1in/1in/0.5in/0.5in -> 54/54/72/72/36/36 points).
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Columns B:G hold fund codes / names / figures that must stay as literal
# text (fund codes like "002345" have significant leading zeros, and the
# numeric-looking figures are stored as text in this workbook's format) -
# force text formatting before writing so Excel doesn't coerce them to
# numbers.
$ws.Range("B1:G21").NumberFormat = "@"

$ws.Cells.Item(1,2).Value = '基金代码'
$ws.Cells.Item(1,3).Value = '基金名称'
$ws.Cells.Item(1,4).Value = '基金规模'
$ws.Cells.Item(1,5).Value = '股票总仓位'
$ws.Cells.Item(1,6).Value = '仓位占比'
$ws.Cells.Item(1,7).Value = '持有市值(亿元)'
$ws.Cells.Item(1,8).Value = '仓位排名'
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = '340007'
$ws.Cells.Item(2,3).Value = '兴全社会责任混合'
$ws.Cells.Item(2,4).Value = '45.11'
$ws.Cells.Item(2,5).Value = '93.56'
$ws.Cells.Item(2,6).Value = '3.18'
$ws.Cells.Item(2,7).Value = '1.4345'
$ws.Cells.Item(2,8).Value = 10
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = '002345'
$ws.Cells.Item(3,3).Value = '华夏高端制造灵活配置混合A'
$ws.Cells.Item(3,4).Value = '16.65'
$ws.Cells.Item(3,5).Value = '93.68'
$ws.Cells.Item(3,6).Value = '7.09'
$ws.Cells.Item(3,7).Value = '1.1805'
$ws.Cells.Item(3,8).Value = 3
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = '070013'
$ws.Cells.Item(4,3).Value = '嘉实研究精选混合A'
$ws.Cells.Item(4,4).Value = '12.63'
$ws.Cells.Item(4,5).Value = '90.55'
$ws.Cells.Item(4,6).Value = '2.48'
$ws.Cells.Item(4,7).Value = '0.3132'
$ws.Cells.Item(4,8).Value = 8
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = '004640'
$ws.Cells.Item(5,3).Value = '华夏节能环保股票A'
$ws.Cells.Item(5,4).Value = '6.25'
$ws.Cells.Item(5,5).Value = '93.24'
$ws.Cells.Item(5,6).Value = '4.71'
$ws.Cells.Item(5,7).Value = '0.2944'
$ws.Cells.Item(5,8).Value = 8
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = '012703'
$ws.Cells.Item(6,3).Value = '华夏核心成长混合A'
$ws.Cells.Item(6,4).Value = '4.93'
$ws.Cells.Item(6,5).Value = '92.86'
$ws.Cells.Item(6,6).Value = '5.06'
$ws.Cells.Item(6,7).Value = '0.2495'
$ws.Cells.Item(6,8).Value = 9
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = '015229'
$ws.Cells.Item(7,3).Value = '华夏低碳经济一年持有混合A'
$ws.Cells.Item(7,4).Value = '2.55'
$ws.Cells.Item(7,5).Value = '86.07'
$ws.Cells.Item(7,6).Value = '4.68'
$ws.Cells.Item(7,7).Value = '0.1193'
$ws.Cells.Item(7,8).Value = 9
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = '016590'
$ws.Cells.Item(8,3).Value = '富国汽车智选混合A'
$ws.Cells.Item(8,4).Value = '5.06'
$ws.Cells.Item(8,5).Value = '49.95'
$ws.Cells.Item(8,6).Value = '2.22'
$ws.Cells.Item(8,7).Value = '0.1123'
$ws.Cells.Item(8,8).Value = 9
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = '013108'
$ws.Cells.Item(9,3).Value = '华夏先进制造龙头混合C'
$ws.Cells.Item(9,4).Value = '2.50'
$ws.Cells.Item(9,5).Value = '90.41'
$ws.Cells.Item(9,6).Value = '4.39'
$ws.Cells.Item(9,7).Value = '0.1098'
$ws.Cells.Item(9,8).Value = 8
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = '013107'
$ws.Cells.Item(10,3).Value = '华夏先进制造龙头混合A'
$ws.Cells.Item(10,4).Value = '1.82'
$ws.Cells.Item(10,5).Value = '90.41'
$ws.Cells.Item(10,6).Value = '4.39'
$ws.Cells.Item(10,7).Value = '0.0799'
$ws.Cells.Item(10,8).Value = 8
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = '673141'
$ws.Cells.Item(11,3).Value = '西部利得景程灵活配置混合A'
$ws.Cells.Item(11,4).Value = '1.42'
$ws.Cells.Item(11,5).Value = '82.89'
$ws.Cells.Item(11,6).Value = '3.56'
$ws.Cells.Item(11,7).Value = '0.0506'
$ws.Cells.Item(11,8).Value = 9
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = '015230'
$ws.Cells.Item(12,3).Value = '华夏低碳经济一年持有混合C'
$ws.Cells.Item(12,4).Value = '1.05'
$ws.Cells.Item(12,5).Value = '86.07'
$ws.Cells.Item(12,6).Value = '4.68'
$ws.Cells.Item(12,7).Value = '0.0491'
$ws.Cells.Item(12,8).Value = 9
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = '015058'
$ws.Cells.Item(13,3).Value = '华夏高端制造灵活配置混合C'
$ws.Cells.Item(13,4).Value = '0.68'
$ws.Cells.Item(13,5).Value = '93.68'
$ws.Cells.Item(13,6).Value = '7.09'
$ws.Cells.Item(13,7).Value = '0.0482'
$ws.Cells.Item(13,8).Value = 3
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = '012710'
$ws.Cells.Item(14,3).Value = '华夏核心成长混合C'
$ws.Cells.Item(14,4).Value = '0.60'
$ws.Cells.Item(14,5).Value = '92.86'
$ws.Cells.Item(14,6).Value = '5.06'
$ws.Cells.Item(14,7).Value = '0.0304'
$ws.Cells.Item(14,8).Value = 9
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = '673143'
$ws.Cells.Item(15,3).Value = '西部利得景程灵活配置混合C'
$ws.Cells.Item(15,4).Value = '0.62'
$ws.Cells.Item(15,5).Value = '82.89'
$ws.Cells.Item(15,6).Value = '3.56'
$ws.Cells.Item(15,7).Value = '0.0221'
$ws.Cells.Item(15,8).Value = 9
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = '015060'
$ws.Cells.Item(16,3).Value = '华夏节能环保股票C'
$ws.Cells.Item(16,4).Value = '0.38'
$ws.Cells.Item(16,5).Value = '93.24'
$ws.Cells.Item(16,6).Value = '4.71'
$ws.Cells.Item(16,7).Value = '0.0179'
$ws.Cells.Item(16,8).Value = 8
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = '002604'
$ws.Cells.Item(17,3).Value = '华夏新起点灵活配置混合A'
$ws.Cells.Item(17,4).Value = '0.56'
$ws.Cells.Item(17,5).Value = '30.09'
$ws.Cells.Item(17,6).Value = '2.97'
$ws.Cells.Item(17,7).Value = '0.0166'
$ws.Cells.Item(17,8).Value = 4
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = '004536'
$ws.Cells.Item(18,3).Value = '嘉实中小企业量化活力灵活配置混合'
$ws.Cells.Item(18,4).Value = '0.22'
$ws.Cells.Item(18,5).Value = '93.77'
$ws.Cells.Item(18,6).Value = '3.91'
$ws.Cells.Item(18,7).Value = '0.0086'
$ws.Cells.Item(18,8).Value = 8
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = '016591'
$ws.Cells.Item(19,3).Value = '富国汽车智选混合C'
$ws.Cells.Item(19,4).Value = '0.38'
$ws.Cells.Item(19,5).Value = '49.95'
$ws.Cells.Item(19,6).Value = '2.22'
$ws.Cells.Item(19,7).Value = '0.0084'
$ws.Cells.Item(19,8).Value = 9
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = '008213'
$ws.Cells.Item(20,3).Value = '华夏新起点灵活配置混合C'
$ws.Cells.Item(20,4).Value = '0.01'
$ws.Cells.Item(20,5).Value = '30.09'
$ws.Cells.Item(20,6).Value = '2.97'
$ws.Cells.Item(20,7).Value = '0.0003'
$ws.Cells.Item(20,8).Value = 4
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = '960025'
$ws.Cells.Item(21,3).Value = '嘉实研究精选混合H'
$ws.Cells.Item(21,4).Value = '0.00'
$ws.Cells.Item(21,5).Value = '90.55'
$ws.Cells.Item(21,6).Value = '2.48'
$ws.Cells.Item(21,7).Value = 0
$ws.Cells.Item(21,8).Value = 8

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new leading data row for
#    2022-Q4 and shift the existing quarters' rows down by one.
#    Final layout (row -> A, B, C, D):
#      2 -> 0, 2022-Q4, 20, 4.15
#      3 -> 1, 2022-Q3, 13, 2.11
#      4 -> 2, 2022-Q2, 16, 3.34
#      5 -> 3, 2022-Q1,  1, 0.78
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 20
$total.Cells.Item(2,4).Value = 4.15

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2022-Q3"
$total.Cells.Item(3,3).Value = 13
$total.Cells.Item(3,4).Value = 2.11

$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2022-Q2"
$total.Cells.Item(4,3).Value = 16
$total.Cells.Item(4,4).Value = 3.34

$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(5,2).Value = "2022-Q1"
$total.Cells.Item(5,3).Value = 1
$total.Cells.Item(5,4).Value = 0.78

# ---------------------------------------------------------------------------
# 3) Restore "2022-Q1" (last sheet) as the active/selected sheet, matching
#    the original workbook's selection state.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Activate()
